# Apply targeted numeric corrections to Sheet1 per the commit "to clear and finish".
# Columns: A=Pregnancies, B=Glucose, C=BMI, D=DiabetesPedigreeFunction, E=Age

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 13.35
$ws.Range("D3").Value = 1.2
$ws.Range("C21").Value = 50.55
$ws.Range("D32").Value = 1.2
$ws.Range("C44").Value = 13.35
$ws.Range("D51").Value = 1.2
$ws.Range("C64").Value = 50.55
$ws.Range("D84").Value = 1.2
$ws.Range("B99").Value = 37.125
$ws.Range("E135").Value = 66.5
$ws.Range("C140").Value = 50.55
$ws.Range("B154").Value = 37.125
$ws.Range("E155").Value = 66.5
$ws.Range("C156").Value = 13.35
$ws.Range("C157").Value = 13.35
$ws.Range("C171").Value = 13.35
$ws.Range("B177").Value = 37.125
$ws.Range("A185").Value = 13.5
$ws.Range("C188").Value = 50.55
$ws.Range("D200").Value = 1.2
$ws.Range("B218").Value = 37.125
$ws.Range("C219").Value = 13.35
$ws.Range("D219").Value = 1.2
$ws.Range("D229").Value = 1.2
$ws.Range("C230").Value = 50.55
$ws.Range("D230").Value = 1.2
$ws.Range("E244").Value = 66.5
$ws.Range("D264").Value = 1.2
$ws.Range("D274").Value = 1.2
$ws.Range("A275").Value = 13.5
$ws.Range("D301").Value = 1.2
$ws.Range("E313").Value = 66.5
$ws.Range("D318").Value = 1.2
$ws.Range("D327").Value = 1.2
$ws.Range("D352").Value = 1.2
$ws.Range("C391").Value = 50.55
$ws.Range("C400").Value = 13.35
$ws.Range("E400").Value = 66.5
$ws.Range("D403").Value = 1.2
$ws.Range("E412").Value = 66.5
$ws.Range("A473").Value = 13.5
$ws.Range("D493").Value = 1.2
$ws.Range("B501").Value = 37.125
$ws.Range("E509").Value = 66.5
$ws.Range("D529").Value = 1.2
$ws.Range("A536").Value = 13.5
$ws.Range("D573").Value = 1.2
$ws.Range("D582").Value = 1.2
$ws.Range("E588").Value = 66.5
$ws.Range("D596").Value = 1.2
$ws.Range("D597").Value = 1.2
$ws.Range("D604").Value = 1.2
